$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("M20").Formula = "'1"
$ws.Range("M21").Formula = "'1"
$ws.Range("M22").Formula = "'2"
